$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create a new sheet "05-03-2022" by copying the existing blank
#    "Daily Attendance Template" sheet, and place it right before the
#    template sheet (so order becomes Key, 05-03-2022, Daily Attendance
#    Template). This is effectively "today's sheet".
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("Daily Attendance Template")
$template.Copy($template, $null)
$today = $wb.Worksheets.Item("Daily Attendance Template (2)")
$today.Name = "05-03-2022"

# The template carries sheet protection - unlock the new copy so data can
# be written into it (matches the resulting file, which has no
# sheetProtection element on the filled-in sheet).
$today.Unprotect()

# ---------------------------------------------------------------------------
# 2. Populate the new sheet with a day's worth of attendance data.
# ---------------------------------------------------------------------------

# Row 2: Bunk 2 / Staff Member 4
$today.Range("A2").Value = "Bunk 2"
$today.Range("B2").Value = "Staff Member 4"
$today.Range("D2").Value = "10:02 PM"
$today.Range("E2").Value2 = 0.92013888888888884
$today.Range("E2").NumberFormat = "h:mm AM/PM"
$today.Range("I2").Value = "1:00 AM"

# Row 3: Bunk 1 / Staff Member 1
$today.Range("A3").Value = "Bunk 1"
$today.Range("B3").Value = "Staff Member 1"
$today.Range("D3").Value = "10:03 PM"
$today.Range("E3").Value = "10:40 PM"
$today.Range("E3").Interior.Color = 13492663
$today.Range("I3").Value = "1:00 AM"

# Row 4: Visitor / Visitor 1
$today.Range("A4").Value = "Visitor"
$today.Range("B4").Value = "Visitor 1"
$today.Range("D4").Value = "10:40 PM (visitor)"
$today.Range("E4").Value2 = 0.9194444444444444
$today.Range("E4").NumberFormat = "h:mm AM/PM"
$today.Range("I4").Value = "5:00 PM"

# Summary counters on the right hand side of the sheet.
$today.Range("I6").Value = 2
$today.Range("I7").Value = 2
$today.Range("I8").Value = 0
$today.Range("I10").Value = 0

$today.Range("I9").Select()

# ---------------------------------------------------------------------------
# 3. Mark Staff Member 1 as "On Time" on the Key sheet.
# ---------------------------------------------------------------------------
$key = $wb.Worksheets.Item("Key")
$key.Range("D2").Value = 1
